# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 349
$wsOff.Range("C3").Value = 226
$wsOff.Range("D3").Value = 150
$wsOff.Range("E3").Value = 61
$wsOff.Range("F3").Value = 7

# DEF sheet - row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 473
$wsDef.Range("C3").Value = 360
$wsDef.Range("D3").Value = 111
$wsDef.Range("E3").Value = 55
